$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout update -------------------------------------------------
# The template gains one new column "Nhân viên nhận" right before the
# existing "Nhân viên phát" column (old column K), and four new columns
# "Dịch vụ cộng thêm 2".."Dịch vụ cộng thêm 5" right after the existing
# "Dịch vụ cộng thêm 1" column. Everything else shifts right accordingly.
# Inserting whole columns re-uses the formatting of the following column,
# which matches the style pattern seen in the target file.

# Insert 1 column at K (before "Nhân viên phát") -> becomes "Nhân viên nhận"
$ws.Columns("K:K").Insert() | Out-Null

# After the above insert, "Dịch vụ cộng thêm 1" now sits in column O, so
# insert 4 columns right after it, at P:S.
$ws.Columns("P:S").Insert() | Out-Null

# --- Fill in the headers for the newly inserted columns -------------------
# Set the "Dịch vụ cộng thêm 2..5" values before "Nhân viên nhận" so that
# shared-string table indices come out in the same order as the source file.
$ws.Range("P1").Value = "Dịch vụ cộng thêm 2"
$ws.Range("Q1").Value = "Dịch vụ cộng thêm 3"
$ws.Range("R1").Value = "Dịch vụ cộng thêm 4"
$ws.Range("S1").Value = "Dịch vụ cộng thêm 5"
$ws.Range("K1").Value = "Nhân viên nhận"

# --- Update the current selection on the sheet -----------------------------
$ws.Range("K1").Select() | Out-Null
